# Auto-generated edit script: updates computed profit/price columns
# (H..N: currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
#  LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) across several
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed
# market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6025
$ws.Range("J74").Value = 6000
$ws.Range("L74").Value = 6000
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 6025
$ws.Range("J77").Value = 6000
$ws.Range("L77").Value = 30000
$ws.Range("N77").Value = -39360
$ws.Range("H98").Value = 3599.762
$ws.Range("I98").Value = 2878.7878
$ws.Range("J98").Value = 6243.3335
$ws.Range("K98").Value = 2878.7878
$ws.Range("L98").Value = 6243.3335
$ws.Range("M98").Value = -1380.7878
$ws.Range("N98").Value = -9239.333500000001
$ws.Range("H122").Value = 3599.762
$ws.Range("I122").Value = 2878.7878
$ws.Range("J122").Value = 6243.3335
$ws.Range("K122").Value = 8636.3634
$ws.Range("L122").Value = 18730.0005
$ws.Range("M122").Value = -6186.3634
$ws.Range("N122").Value = -23630.0005
$ws.Range("H137").Value = 4266.8423
$ws.Range("I137").Value = 2057.8
$ws.Range("J137").Value = 12550.75
$ws.Range("K137").Value = 6173.400000000001
$ws.Range("L137").Value = 37652.25
$ws.Range("M137").Value = -3623.400000000001
$ws.Range("N137").Value = -42752.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H32").Value = 3219.4824
$ws.Range("I32").Value = 2670
$ws.Range("K32").Value = 2670
$ws.Range("M32").Value = -2383
$ws.Range("H61").Value = 20002584
$ws.Range("I61").Value = 21741568
$ws.Range("J61").Value = 4287.5
$ws.Range("K61").Value = 21741568
$ws.Range("L61").Value = 4287.5
$ws.Range("M61").Value = -21741356
$ws.Range("N61").Value = -4711.5
$ws.Range("H112").Value = 13443.5
$ws.Range("J112").Value = 13443.5
$ws.Range("L112").Value = 13443.5
$ws.Range("N112").Value = -16397.5
$ws.Range("H132").Value = 31304664
$ws.Range("I132").Value = 13027.28
$ws.Range("J132").Value = 143060510
$ws.Range("K132").Value = 39081.84
$ws.Range("L132").Value = 429181530
$ws.Range("M132").Value = -36551.84
$ws.Range("N132").Value = -429186590
$ws.Range("H136").Value = 20002584
$ws.Range("I136").Value = 21741568
$ws.Range("J136").Value = 4287.5
$ws.Range("K136").Value = 65224704
$ws.Range("L136").Value = 12862.5
$ws.Range("M136").Value = -65222154
$ws.Range("N136").Value = -17962.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1828.5883
$ws.Range("I94").Value = 1507.5
$ws.Range("J94").Value = 2599.2
$ws.Range("K94").Value = 1507.5
$ws.Range("L94").Value = 2599.2
$ws.Range("M94").Value = -1056.5
$ws.Range("N94").Value = -3501.2
$ws.Range("H99").Value = 4178.706
$ws.Range("I99").Value = 2751.5833
$ws.Range("J99").Value = 7603.8
$ws.Range("K99").Value = 2751.5833
$ws.Range("L99").Value = 7603.8
$ws.Range("M99").Value = -1253.5833
$ws.Range("N99").Value = -10599.8
$ws.Range("H107").Value = 3733.1667
$ws.Range("I107").Value = 2815
$ws.Range("J107").Value = 4651.3335
$ws.Range("K107").Value = 2815
$ws.Range("L107").Value = 4651.3335
$ws.Range("M107").Value = -895
$ws.Range("N107").Value = -8491.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7146605.5
$ws.Range("I31").Value = 2377.9565
$ws.Range("J31").Value = 20839708
$ws.Range("K31").Value = 2377.9565
$ws.Range("L31").Value = 20839708
$ws.Range("M31").Value = -2082.9565
$ws.Range("N31").Value = -20840298
$ws.Range("H34").Value = 7146605.5
$ws.Range("I34").Value = 2377.9565
$ws.Range("J34").Value = 20839708
$ws.Range("K34").Value = 2377.9565
$ws.Range("L34").Value = 20839708
$ws.Range("M34").Value = -2175.9565
$ws.Range("N34").Value = -20840112
$ws.Range("H58").Value = 1956.9333
$ws.Range("I58").Value = 1826.2727
$ws.Range("J58").Value = 2316.25
$ws.Range("K58").Value = 1826.2727
$ws.Range("L58").Value = 2316.25
$ws.Range("M58").Value = -1623.2727
$ws.Range("N58").Value = -2722.25
$ws.Range("H99").Value = 8079.7915
$ws.Range("I99").Value = 5036.3
$ws.Range("K99").Value = 5036.3
$ws.Range("M99").Value = -3538.3
$ws.Range("H100").Value = 79999.75
$ws.Range("J100").Value = 79999.75
$ws.Range("L100").Value = 79999.75
$ws.Range("N100").Value = -82163.75
$ws.Range("H126").Value = 8079.7915
$ws.Range("I126").Value = 5036.3
$ws.Range("K126").Value = 15108.9
$ws.Range("M126").Value = -12638.9
$ws.Range("H136").Value = 1956.9333
$ws.Range("I136").Value = 1826.2727
$ws.Range("J136").Value = 2316.25
$ws.Range("K136").Value = 5478.8181
$ws.Range("L136").Value = 6948.75
$ws.Range("M136").Value = -2928.8181
$ws.Range("N136").Value = -12048.75
$ws.Range("H141").Value = 115578.625
$ws.Range("I141").Value = 55000
$ws.Range("J141").Value = 124232.71
$ws.Range("K141").Value = 55000
$ws.Range("L141").Value = 124232.71
$ws.Range("M141").Value = -49820
$ws.Range("N141").Value = -134592.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 710.6429000000001
$ws.Range("I12").Value = 257.66666
$ws.Range("J12").Value = 1050.375
$ws.Range("K12").Value = 772.9999799999999
$ws.Range("L12").Value = 3151.125
$ws.Range("M12").Value = -599.9999799999999
$ws.Range("N12").Value = -3497.125
$ws.Range("H107").Value = 1704.1111
$ws.Range("J107").Value = 1850.5
$ws.Range("L107").Value = 5551.5
$ws.Range("N107").Value = -9391.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 137974.53
$ws.Range("I70").Value = 252902.38
$ws.Range("J70").Value = 6628.4287
$ws.Range("K70").Value = 252902.38
$ws.Range("L70").Value = 6628.4287
$ws.Range("M70").Value = -252632.38
$ws.Range("N70").Value = -7168.4287
$ws.Range("H73").Value = 137974.53
$ws.Range("I73").Value = 252902.38
$ws.Range("J73").Value = 6628.4287
$ws.Range("K73").Value = 252902.38
$ws.Range("L73").Value = 6628.4287
$ws.Range("M73").Value = -251966.38
$ws.Range("N73").Value = -8500.4287
$ws.Range("H80").Value = 13150.5
$ws.Range("I80").Value = 15999.667
$ws.Range("J80").Value = 10301.333
$ws.Range("K80").Value = 15999.667
$ws.Range("L80").Value = 10301.333
$ws.Range("M80").Value = -15001.667
$ws.Range("N80").Value = -12297.333
$ws.Range("H83").Value = 13150.5
$ws.Range("I83").Value = 15999.667
$ws.Range("J83").Value = 10301.333
$ws.Range("K83").Value = 79998.33499999999
$ws.Range("L83").Value = 51506.665
$ws.Range("M83").Value = -75006.33499999999
$ws.Range("N83").Value = -61490.665
$ws.Range("H107").Value = 687.1667
$ws.Range("J107").Value = 612.7143
$ws.Range("L107").Value = 612.7143
$ws.Range("N107").Value = -4452.7143
$ws.Range("H132").Value = 5113.7856
$ws.Range("I132").Value = 5217.795
$ws.Range("J132").Value = 3761.6667
$ws.Range("K132").Value = 15653.385
$ws.Range("L132").Value = 11285.0001
$ws.Range("M132").Value = -13123.385
$ws.Range("N132").Value = -16345.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 86666
$ws.Range("I81").Value = 84999.5
$ws.Range("J81").Value = 89999
$ws.Range("K81").Value = 84999.5
$ws.Range("L81").Value = 89999
$ws.Range("M81").Value = -84001.5
$ws.Range("N81").Value = -91995
$ws.Range("H84").Value = 86666
$ws.Range("I84").Value = 84999.5
$ws.Range("J84").Value = 89999
$ws.Range("K84").Value = 254998.5
$ws.Range("L84").Value = 269997
$ws.Range("M84").Value = -250006.5
$ws.Range("N84").Value = -279981
$ws.Range("H110").Value = 14457.6
$ws.Range("J110").Value = 14457.6
$ws.Range("L110").Value = 14457.6
$ws.Range("N110").Value = -22637.6
$ws.Range("H122").Value = 3293472.5
$ws.Range("I122").Value = 3834.9
$ws.Range("K122").Value = 11504.7
$ws.Range("M122").Value = -9054.700000000001
$ws.Range("H131").Value = 73433.39999999999
$ws.Range("J131").Value = 87290.336
$ws.Range("L131").Value = 87290.336
$ws.Range("N131").Value = -97370.336
$ws.Range("H132").Value = 2668.5908
$ws.Range("J132").Value = 3457.1428
$ws.Range("L132").Value = 10371.4284
$ws.Range("N132").Value = -15431.4284
$ws.Range("H136").Value = 2005479.2
$ws.Range("I136").Value = 2225199
$ws.Range("J136").Value = 28000
$ws.Range("K136").Value = 6675597
$ws.Range("L136").Value = 84000
$ws.Range("M136").Value = -6673047
$ws.Range("N136").Value = -89100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 32998
$ws.Range("J70").Value = 32998
$ws.Range("L70").Value = 32998
$ws.Range("N70").Value = -33628
$ws.Range("H73").Value = 32998
$ws.Range("J73").Value = 32998
$ws.Range("L73").Value = 32998
$ws.Range("N73").Value = -35182
$ws.Range("H76").Value = 35000
$ws.Range("I76").Value = 35000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 35000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -34685
$ws.Range("N76").Value = $null
$ws.Range("H79").Value = 35000
$ws.Range("I79").Value = 35000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 35000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -33908
$ws.Range("N79").Value = $null
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("H122").Value = 25006300
$ws.Range("I122").Value = 3579.6
$ws.Range("J122").Value = 66677504
$ws.Range("K122").Value = 10738.8
$ws.Range("L122").Value = 200032512
$ws.Range("M122").Value = -8288.799999999999
$ws.Range("N122").Value = -200037412
$ws.Range("H124").Value = 3370666
$ws.Range("J124").Value = 3370666
$ws.Range("L124").Value = 3370666
$ws.Range("N124").Value = -3380486
$ws.Range("H131").Value = 87331.336
$ws.Range("J131").Value = 87331.336
$ws.Range("L131").Value = 87331.336
$ws.Range("N131").Value = -97411.336
$ws.Range("H132").Value = 2270.5715
$ws.Range("I132").Value = 2162.5454
$ws.Range("K132").Value = 6487.6362
$ws.Range("M132").Value = -3957.6362
$ws.Range("H136").Value = 2724.652
$ws.Range("I136").Value = 2342.8333
$ws.Range("J136").Value = 4099.2
$ws.Range("K136").Value = 7028.499899999999
$ws.Range("L136").Value = 12297.6
$ws.Range("M136").Value = -4478.499899999999
$ws.Range("N136").Value = -17397.6
